# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
# Commit: Updated cryptos list on Sat May 18 11:51:36 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.243.18"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.145.31"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.31"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.72"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.139.20"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  +2.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.484"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.48"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "3.659.89"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "67.219.51"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "3.139.74"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.18"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.88"
$ws.Range("E21").Value = "  +5.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.34"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.27"
$ws.Range("E25").Value = "  +4.09%  "
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.06"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").Value = "0.0₂01000"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.991"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.41"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.313"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "2.851.02"
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "385.23"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.75"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.94"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("E51").Value = "  -0.20%  "
